# Updates the cryptos price/volume snapshot table (Thu Apr 20 17:47:43 UTC 2023
# GitHub Actions refresh): Price (column D) and Volume(1h) (column E) values are
# refreshed for rows 2-51, and two coin pairs (row 39/40 and row 49/50) swap
# ranking positions, so their Coin/Link/Price/Volume cells are updated in place.
# All Price/Volume cells are plain text in the workbook (numbers such as
# '28.829.61' or '1.013' are not valid numeric values and the percentages keep
# their surrounding padding), so every write below is forced to text with a
# leading apostrophe to stop Excel from reinterpreting it as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.829.61"
$ws.Range("E2").Value = "'  -2.53%  "
$ws.Range("D3").Value = "'1.964.19"
$ws.Range("E3").Value = "'  -1.82%  "
$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "'  -0.12%  "
$ws.Range("D5").Value = "'323.68"
$ws.Range("E5").Value = "'  -2.06%  "
$ws.Range("D6").Value = "'1.013"
$ws.Range("E6").Value = "'  -0.05%  "
$ws.Range("E7").Value = "'  -4.47%  "
$ws.Range("D8").Value = "'0.4037"
$ws.Range("D9").Value = "'53.96"
$ws.Range("E9").Value = "'  -0.48%  "
$ws.Range("D10").Value = "'0.08472"
$ws.Range("E10").Value = "'  -5.41%  "
$ws.Range("E11").Value = "'  -5.20%  "
$ws.Range("D12").Value = "'22.43"
$ws.Range("E12").Value = "'  -4.02%  "
$ws.Range("D13").Value = "'1.950.83"
$ws.Range("E13").Value = "'  -4.57%  "
$ws.Range("D14").Value = "'7.656"
$ws.Range("E14").Value = "'  -4.97%  "
$ws.Range("D15").Value = "'6.239"
$ws.Range("E15").Value = "'  -3.92%  "
$ws.Range("D16").Value = "'1.014"
$ws.Range("E16").Value = "'  -0.01%  "
$ws.Range("D17").Value = "'89.83"
$ws.Range("E17").Value = "'  -4.81%  "
$ws.Range("D18").Value = "'0.00001069"
$ws.Range("E18").Value = "'  -3.97%  "
$ws.Range("D19").Value = "'0.06597"
$ws.Range("E19").Value = "'  -1.15%  "
$ws.Range("D20").Value = "'18.67"
$ws.Range("E20").Value = "'  -5.55%  "
$ws.Range("D21").Value = "'1.011"
$ws.Range("E21").Value = "'  -0.15%  "
$ws.Range("D22").Value = "'5.788"
$ws.Range("E22").Value = "'  -2.89%  "
$ws.Range("D23").Value = "'28.831.74"
$ws.Range("E23").Value = "'  -2.58%  "
$ws.Range("D24").Value = "'11.54"
$ws.Range("E24").Value = "'  -3.71%  "
$ws.Range("D25").Value = "'2.291"
$ws.Range("E25").Value = "'  +0.53%  "
$ws.Range("D26").Value = "'2.202.21"
$ws.Range("E26").Value = "'  -3.38%  "
$ws.Range("D27").Value = "'154.59"
$ws.Range("E27").Value = "'  -3.12%  "
$ws.Range("D28").Value = "'20.22"
$ws.Range("E28").Value = "'  -2.30%  "
$ws.Range("D29").Value = "'5.959"
$ws.Range("E29").Value = "'  -5.95%  "
$ws.Range("D30").Value = "'2.158"
$ws.Range("E30").Value = "'  -6.36%  "
$ws.Range("D31").Value = "'124.13"
$ws.Range("E31").Value = "'  -3.36%  "
$ws.Range("E32").Value = "'  -4.62%  "
$ws.Range("D33").Value = "'0.09614"
$ws.Range("E33").Value = "'  -3.38%  "
$ws.Range("D34").Value = "'1.460"
$ws.Range("E34").Value = "'  -6.79%  "
$ws.Range("D35").Value = "'5.683"
$ws.Range("E35").Value = "'  -2.78%  "
$ws.Range("D36").Value = "'3.687"
$ws.Range("E36").Value = "'  -2.59%  "
$ws.Range("D37").Value = "'0.02360"
$ws.Range("E37").Value = "'  -4.23%  "
$ws.Range("D38").Value = "'1.267"
$ws.Range("E38").Value = "'  -3.19%  "
$ws.Range("B39").Value = "'Hedera"
$ws.Range("C39").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06210"
$ws.Range("E39").Value = "'  -2.32%  "
$ws.Range("B40").Value = "'FraxShare"
$ws.Range("C40").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'8.767"
$ws.Range("E40").Value = "'  -7.10%  "
$ws.Range("D41").Value = "'0.6240"
$ws.Range("E41").Value = "'  -5.05%  "
$ws.Range("D42").Value = "'11.12"
$ws.Range("E42").Value = "'  -4.78%  "
$ws.Range("E43").Value = "'  -0.26%  "
$ws.Range("E44").Value = "'  -6.59%  "
$ws.Range("D45").Value = "'1.344"
$ws.Range("E45").Value = "'  +2.43%  "
$ws.Range("D46").Value = "'0.5965"
$ws.Range("E46").Value = "'  -5.83%  "
$ws.Range("D47").Value = "'12.95"
$ws.Range("E47").Value = "'  -4.22%  "
$ws.Range("D48").Value = "'2.077"
$ws.Range("E48").Value = "'  -5.86%  "
$ws.Range("B49").Value = "'PancakeSwap"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D49").Value = "'3.424"
$ws.Range("E49").Value = "'  -2.84%  "
$ws.Range("B50").Value = "'BabyDogeCoin"
$ws.Range("C50").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.00000000334"
$ws.Range("E50").Value = "'  -2.03%  "
$ws.Range("D51").Value = "'0.06837"
$ws.Range("E51").Value = "'  -2.18%  "
